$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The edit swaps the per-occurrence data between row 10 and row 11 (the
# observation records), while a handful of shared columns (C, D, L, N, S,
# T, U, V, W, Y, Z, AA, AB, AD, AE, AF, AG, AT, AW, AX, AY) stay the same.
# Capture the two rows' current values for the columns that move, then
# write them back swapped.

$cols = @("A","B","E","F","G","H","I","J","K","P","Q","R","AI")

$row10 = @{}
$row11 = @{}
foreach ($col in $cols) {
    $row10[$col] = $ws.Range("${col}10").Value2
    $row11[$col] = $ws.Range("${col}11").Value2
}

foreach ($col in $cols) {
    $ws.Range("${col}10").Value2 = $row11[$col]
    $ws.Range("${col}11").Value2 = $row10[$col]
}
